$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tuesday")

$ws.Range("B2").Value = "tuesday pronunciation american"
$ws.Range("C2").Value = "tuesday"

$ws.Range("B3").Value = "baby pink color"
$ws.Range("C3").Value = "Baby"
